$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the vendor entries "POP" (A13) and "zuluCrypt" (A14) while keeping
# their existing formatting/style.
$ws.Range("A13:A14").ClearContents()

# Move the active selection to A5.
$ws.Range("A5").Select()
